$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Materials")

# Update the scientificNameAuthorship mapping value (row 2) from
# ${summary.Author} to ${summary.authority}
$ws.Range("BB2").Value = "`${summary.authority}"

# Remove the Taxon_Local_ID column (A) and the suborder/infraorder/superfamily
# columns (AR/AS/AT) entirely. Delete from right to left so earlier column
# letters stay valid while later ones are removed.
$ws.Range("AT1").EntireColumn.Delete()
$ws.Range("AS1").EntireColumn.Delete()
$ws.Range("AR1").EntireColumn.Delete()
$ws.Range("A1").EntireColumn.Delete()
